# Add 2022-Q3 data
# 1) Insert a new row into the "总计" (summary) sheet with the new quarter's totals,
#    shifting the existing rows down by one.
# 2) Insert a brand-new worksheet named "2022-Q3" right after "总计" holding the
#    per-fund holdings detail for the new quarter, pushing the older quarter
#    sheets one position to the right.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Part 1: "总计" sheet - insert new top data row (2022-Q3: 18 holdings, 3.4 亿元)
# and push every existing data row down by one, bumping its running index (col A).
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)
$summary.Rows.Item(2).Insert()
# Insert() copies the header row's formatting down onto the new blank row;
# strip that back off the B:D cells so they land with no explicit style,
# matching every other data row in this column range.
$summary.Range("B2:D2").ClearFormats()

$summaryRows = @(
    @("2022-Q3", 18, 3.4),
    @("2022-Q2", 17, 4.41),
    @("2022-Q1", 14, 3.38),
    @("2021-Q4", 13, 3.64),
    @("2021-Q3", 20, 3.5),
    @("2021-Q2", 14, 2.88),
    @("2021-Q1", 12, 2.33),
    @("2020-Q4", 21, 5.26)
)

$r = 2
foreach ($row in $summaryRows) {
    $idxCell = $summary.Cells.Item($r, 1)
    $idxCell.Value = $r - 2
    $idxCell.Font.Bold = $true
    $idxCell.Borders.LineStyle = 1
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160

    $summary.Cells.Item($r, 2).Value = $row[0]
    $summary.Cells.Item($r, 3).Value = $row[1]
    $summary.Cells.Item($r, 4).Value = $row[2]
    $r++
}

# ---------------------------------------------------------------------------
# Part 2: new "2022-Q3" worksheet with per-fund holdings detail
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($wb.Worksheets.Item(2))
$q3.Name = "2022-Q3"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $cell = $q3.Cells.Item(1, $col)
    $cell.Value = $headers[$col - 2]
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

$rows = @(
    @("513090", "易方达中证香港证券投资主题ETF",          "10.53", "96.33", "15.05", "1.5848", 1),
    @("501050", "华夏沪港通上证50AH优选指数（LOF）A",       "22.49", "91.20", "2.36",  "0.5308", 9),
    @("012943", "广发稳睿六个月持有期混合A",                "20.75", "26.11", "1.28",  "0.2656", 10),
    @("012944", "广发稳睿六个月持有期混合C",                "19.14", "26.11", "1.28",  "0.2450", 10),
    @("005938", "工银瑞信精选金融地产行业混合C",            "5.51",  "89.46", "4.14",  "0.2281", 9),
    @("005937", "工银瑞信精选金融地产行业混合A",            "4.12",  "89.46", "4.14",  "0.1706", 9),
    @("014317", "广发价值领航一年持有混合A",                "2.11",  "91.65", "5.28",  "0.1114", 5),
    @("004497", "前海开源多元策略灵活配置混合C",            "1.68",  "93.04", "3.98",  "0.0669", 9),
    @("011355", "华泰柏瑞港股通时代机遇混合A",              "0.54",  "91.80", "7.79",  "0.0421", 4),
    @("004496", "前海开源多元策略灵活配置混合A",            "0.91",  "93.04", "3.98",  "0.0362", 9),
    @("014318", "广发价值领航一年持有混合C",                "0.58",  "91.65", "5.28",  "0.0306", 5),
    @("003413", "华泰柏瑞新经济沪港深混合",                  "0.42",  "86.45", "5.80",  "0.0244", 4),
    @("460010", "华泰柏瑞亚洲领导企业混合（QDII）",          "0.36",  "93.44", "5.99",  "0.0216", 6),
    @("011356", "华泰柏瑞港股通时代机遇混合C",              "0.24",  "91.80", "7.79",  "0.0187", 4),
    @("001942", "前海开源沪港深汇鑫灵活配置混合A",          "0.17",  "87.24", "4.45",  "0.0076", 9),
    @("006395", "华夏沪港通上证50AH优选指数（LOF）C",       "0.30",  "91.20", "2.36",  "0.0071", 9),
    @("001943", "前海开源沪港深汇鑫灵活配置混合C",          "0.09",  "87.24", "4.45",  "0.0040", 9),
    @("002860", "前海开源沪港深新机遇灵活配置混合",          "0.02",  "42.74", "2.97",  "0.0006", 8)
)

$r = 2
foreach ($row in $rows) {
    $idxCell = $q3.Cells.Item($r, 1)
    $idxCell.Value = $r - 2
    $idxCell.Font.Bold = $true
    $idxCell.Borders.LineStyle = 1
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160

    # Columns B-G are stored as plain text in the source data (fund codes keep
    # leading zeros, percentages/amounts are text, not numbers) - force text
    # number-format first so COM doesn't silently coerce numeric-looking
    # strings (and strip leading zeros from codes like "012943").
    for ($col = 2; $col -le 7; $col++) {
        $q3.Cells.Item($r, $col).NumberFormat = "@"
    }

    $q3.Cells.Item($r, 2).Value = $row[0]
    $q3.Cells.Item($r, 3).Value = $row[1]
    $q3.Cells.Item($r, 4).Value = $row[2]
    $q3.Cells.Item($r, 5).Value = $row[3]
    $q3.Cells.Item($r, 6).Value = $row[4]
    $q3.Cells.Item($r, 7).Value = $row[5]
    $q3.Cells.Item($r, 8).Value = $row[6]
    $r++
}

# Restore the original active sheet ("总计") as the selected tab.
$wb.Worksheets.Item(1).Activate()

Write-Host "2022-Q3 data added"
